# Atualização de bases das ligas, do dia: 14-05-2024 às 20:19
#
# 1) Rows 83 and 84 had their match data swapped (id column A stays 81/82,
#    but all the other columns B..AB - the actual match record - are
#    exchanged between the two rows).
# 2) Four new match rows are appended at the end of the sheet (rows 104-107,
#    ids 102-105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

# --- Step 1: swap the contents of rows 83 and 84 (everything except column A) ---
$row83 = @{}
$row84 = @{}
foreach ($c in $cols) {
    $row83[$c] = $ws.Range($c + "83").Value()
    $row84[$c] = $ws.Range($c + "84").Value()
}
foreach ($c in $cols) {
    $ws.Range($c + "83").Value = $row84[$c]
    $ws.Range($c + "84").Value = $row83[$c]
}

# --- Step 2: append 4 new rows (104-107) re-using row 103's formatting ---
$newRows = @(
    @{ Row=104; A=102; B=7803366; C="Canada Premier League"; D=45422.83333333334; E="York United FC"; F="Valour FC";
       G=3; H=1; I="H"; J=1.95; K=3.4; L=3.25; M=1.8; N=3.5; O=3.75; P=-0.5;
       Q=1.825; R=1.975; S=2.5; T=1.9; U=1.9; V=0.8; W=-1; X=-1; Y=0.825; Z=-1; AA=0.8999999999999999; AB=-1 },
    @{ Row=105; A=103; B=7802941; C="Canada Premier League"; D=45423.625; E="HFX Wanderers"; F="Cavalry FC";
       G=1; H=1; I="D"; J=2.7; K=3.5; L=2.2; M=2.9; N=3.5; O=2.1; P=0.25;
       Q=1.95; R=1.85; S=2.25; T=1.85; U=1.95; V=-1; W=2.5; X=-1; Y=0.475; Z=-0.5; AA=-0.5; AB=0.475 },
    @{ Row=106; A=104; B=7802878; C="Canada Premier League"; D=45423.75; E="Pacific FC CA"; F="Forge FC";
       G=0; H=0; I="D"; J=2.25; K=3.25; L=2.75; M=2.4; N=3.2; O=2.55; P=0;
       Q=1.8; R=2; S=2.5; T=1.95; U=1.85; V=-1; W=2.2; X=-1; Y=0; Z=0; AA=-1; AB=0.8500000000000001 },
    @{ Row=107; A=105; B=7802942; C="Canada Premier League"; D=45424.83333333334; E="Vancouver FC"; F="Atletico Ottawa";
       G=1; H=1; I="D"; J=3.25; K=3.25; L=2; M=2.8; N=3.1; O=2.25; P=0.25;
       Q=1.775; R=2.025; S=2.25; T=1.8; U=2; V=-1; W=2.1; X=-1; Y=0.3875; Z=-0.5; AA=-0.5; AB=0.5 }
)

$allCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

foreach ($rd in $newRows) {
    $r = $rd["Row"]
    $srcRange = "A103:AB103"
    $dstRange = "A" + $r + ":AB" + $r
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4122)

    foreach ($c in $allCols) {
        $ws.Range($c + $r).Value = $rd[$c]
    }
}
